$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old data range first (old range was A1:G2)
$ws.Range("A1:G2").Clear() | Out-Null

# Header row (row 1) - set size increased to 2 (stim2/probe2 columns added)
$ws.Range("A1").Value = "x1"
$ws.Range("B1").Value = "y1"
$ws.Range("C1").Value = "stim1_color"
$ws.Range("D1").Value = "probe1_color"
$ws.Range("E1").Value = "x2"
$ws.Range("F1").Value = "y2"
$ws.Range("G1").Value = "stim2_color"
$ws.Range("H1").Value = "probe2_color"

# Row 2
$ws.Range("A2").Value = -0.25
$ws.Range("B2").Value = -0.25
$ws.Range("C2").Value = "red"
$ws.Range("D2").Value = "blue"
$ws.Range("E2").Value = -0.25
$ws.Range("F2").Value = -0.25
$ws.Range("G2").Value = "orange"
$ws.Range("H2").Value = "orange"

# Row 3 (new row)
$ws.Range("A3").Value = -0.25
$ws.Range("B3").Value = -0.25
$ws.Range("C3").Value = "orange"
$ws.Range("D3").Value = "orange"
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0.25
$ws.Range("G3").Value = "blue"
$ws.Range("H3").Value = "green"

# Update selection to reflect the last active cell after editing
$ws.Range("F7").Select() | Out-Null
